$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (A:D), matching data appended by the telegram bot message sender
$newRows = @(
    @(1, "77C-226.75", "Bottom", "2025-11-22 22:49:08"),
    @(1, "77C-226.75", "Bottom", "2025-11-22 23:40:20"),
    @(35, "77A-247.01", "Top", "2025-11-22 23:40:31")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
